# ProblemSet1 (week1_tues.pptx) — "Version Control System" slide.
#
# The paragraph listing VCS examples currently reads (as a single run):
#   "), Mercury (or Hg), "
# and needs to become four runs:
#   ")"  +  ", "  +  "Mercurial "  +  "(or Hg), "
# i.e. "Mercury" -> "Mercurial", while splitting the surrounding
# punctuation into its own runs (matching how PowerPoint splits a run
# when only part of it is retyped).

$p = $ppt.ActivePresentation

# Slide 12 ("Version Control System"), Shape 2 ("Content Placeholder 2").
$slide = $p.Slides.Item(12)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange

# Locate the run holding "), Mercury (or Hg), " by searching the whole
# shape's text (more robust than walking Paragraphs/Runs collections).
$found = $textRange.Find("), Mercury (or Hg), ", 0)
$runStart = $found.Start

# Old run text, split into the four segments that map onto the new text
# (only the "Mercury " -> "Mercurial " segment actually changes content;
# the others are rewritten in place purely to force the run split):
#   [0,1)   ")"
#   [1,3)   ", "
#   [3,11)  "Mercury "   -> "Mercurial "
#   [11,20) "(or Hg), "
# Apply right-to-left so earlier (left-hand) offsets stay valid while
# later (right-hand) pieces are still being resized.

$seg4 = $textRange.Characters($runStart + 11, 9)
$seg4.Text = "(or Hg), "

$seg3 = $textRange.Characters($runStart + 3, 8)
$seg3.Text = "Mercurial "

$seg2 = $textRange.Characters($runStart + 1, 2)
$seg2.Text = ", "

$seg1 = $textRange.Characters($runStart + 0, 1)
$seg1.Text = ")"
